$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap C16 and D16 values
$ws.Range("C16").Value = "skos:definition@en"
$ws.Range("D16").Value = "skos:altLabel(separator=`",`")"

# Update row 17 with new values
$ws.Range("A17").Value = "nicest-2-subjects:10000"
$ws.Range("B17").Value = "test subject"
$ws.Range("C17").Value = "This is a subject term used for setting up sheet2rdf workflow."
$ws.Range("D17").Value = ""

# Add new row 18
$ws.Range("A18").Value = "nicest-2-subjects:10001"
$ws.Range("B18").Value = "climate"
